$d = $word.ActiveDocument

# --- Edit 1 (route 7 paragraph) ---------------------------------------
# Drop Word's "_GoBack" (last-edit) bookmark at the point the cursor ends
# up at after trimming "also makes a circuit. " out of the sentence --
# right between "...and i" and "t's a little bit...". Doing this before
# any text edit lets the run that carries the leading <w:tab/> split
# cleanly instead of being flattened into plain text.
$r0 = $d.Content
$r0.Find.Execute("to make this route and i")
$goBackRange = $d.Range($r0.End, $r0.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Trim "t also makes a circuit. I" out of the now-split run so
# "...to make this route and it also makes a circuit. It's a little bit..."
# becomes "...to make this route and it's a little bit..."
$r1 = $d.Content
$r1.Find.Execute("t also makes a circuit. I")
$r1.Text = ""

# --- Edit 2 (route 11 paragraph) ---------------------------------------
# Replace the closing clause; temporary bookmarks at the edit boundaries
# force the surrounding text to stay in its own run (matching how the
# unedited lead-in and trailing "10 points..." sentence are split out as
# separate runs) instead of being re-merged into one big run.
$r2 = $d.Content
$r2.Find.Execute("but it comes close enough to the required length to be viable. ")
$editStart = $r2.Start
$editEnd = $r2.End

$d.Bookmarks.Add("ztempStart", $d.Range($editStart, $editStart))
$d.Bookmarks.Add("ztempEnd", $d.Range($editEnd, $editEnd))

$r3 = $d.Range($editStart, $editEnd)
$r3.Text = "but that variation is what we want. We don't want simple routes all the time."

$d.Bookmarks("ztempStart").Delete()
$d.Bookmarks("ztempEnd").Delete()
